$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 and 6 (the extra sample rows), leaving only data row 4
$ws.Rows("5:6").Delete()

# Update row 2 value
$ws.Range("B2").Value = 1

# Update row 4 values to reflect the new computed vector
$ws.Range("F4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("S4").Value = 6
$ws.Range("AB4").Value = [double]"7.582560427911907e-10"
$ws.Range("AC4").Value = 0.002478752176666358
$ws.Range("AD4").Value = 0.04978706836786394
$ws.Range("AJ4").Value = 0.05239569569327056
$ws.Range("CT4").Value = 0
$ws.Range("CU4").Value = [double]"1.879528816539083e-12"
$ws.Range("CV4").Value = [double]"3.775134544279098e-11"
$ws.Range("CX4").Value = 0
$ws.Range("CZ4").Value = 0
$ws.Range("DB4").Value = 0.05239569569327055
$ws.Range("DG4").Value = 0.04978706836786394
$ws.Range("DM4").Value = 0.0001298759447476425
$ws.Range("DW4").Value = 0.002608628083662656
$ws.Range("EL4").Value = 0
$ws.Range("EN4").Value = 0
$ws.Range("FA4").Value = 0
$ws.Range("FJ4").Value = 0.2592331742852723

# Update the label string for row 4 from "test1" to "testmole"
$ws.Range("FV4").Value = "testmole"
